# The source data export regenerated the "PEBCOM" sheet: four closed/duplicate
# cases were removed from the table (rows identified by their original Caso
# values 6085, 5973, 6236 and -499). Every row below each removal shifts up,
# and the sheet's used range shrinks from A1:P76 to A1:P72.
#
# Deleting from the bottom-most row upward keeps the remaining row indices
# stable while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).EntireRow.Delete()
$ws.Rows.Item(71).EntireRow.Delete()
$ws.Rows.Item(59).EntireRow.Delete()
$ws.Rows.Item(34).EntireRow.Delete()
